# Auto-generated Excel COM-interop script to apply Twintania_Profits data updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 266.66666
$ws.Cells.Item(2, 9).Value = 250
$ws.Cells.Item(2, 11).Value = 250
$ws.Cells.Item(2, 13).Value = -137
$ws.Cells.Item(5, 8).Value = 386
$ws.Cells.Item(5, 9).Value = 386
$ws.Cells.Item(5, 11).Value = 386
$ws.Cells.Item(5, 13).Value = -271
$ws.Cells.Item(18, 8).Value = 281.07144
$ws.Cells.Item(18, 9).Value = 281.07144
$ws.Cells.Item(18, 11).Value = 281.07144
$ws.Cells.Item(18, 13).Value = 2.928560000000004
$ws.Cells.Item(32, 8).Value = 4554.2
$ws.Cells.Item(32, 10).Value = 1923.3334
$ws.Cells.Item(32, 12).Value = 1923.3334
$ws.Cells.Item(32, 14).Value = -2575.3334
$ws.Cells.Item(33, 8).Value = 1633.5714
$ws.Cells.Item(33, 9).Value = 776.0625
$ws.Cells.Item(33, 11).Value = 776.0625
$ws.Cells.Item(33, 13).Value = -547.0625
$ws.Cells.Item(40, 8).Value = 1927.8572
$ws.Cells.Item(40, 9).Value = 1750
$ws.Cells.Item(40, 10).Value = 1999
$ws.Cells.Item(40, 11).Value = 1750
$ws.Cells.Item(40, 12).Value = 1999
$ws.Cells.Item(40, 13).Value = -1575
$ws.Cells.Item(40, 14).Value = -2349
$ws.Cells.Item(43, 8).Value = 4678.8
$ws.Cells.Item(43, 9).Value = 2800
$ws.Cells.Item(43, 11).Value = 2800
$ws.Cells.Item(43, 13).Value = -2731
$ws.Cells.Item(88, 8).Value = 10000
$ws.Cells.Item(88, 10).Value = 10000
$ws.Cells.Item(88, 12).Value = 10000
$ws.Cells.Item(88, 14).Value = -10812
$ws.Cells.Item(91, 8).Value = 10000
$ws.Cells.Item(91, 10).Value = 10000
$ws.Cells.Item(91, 12).Value = 10000
$ws.Cells.Item(91, 14).Value = -12808
$ws.Cells.Item(100, 8).Value = 28275.922
$ws.Cells.Item(100, 9).Value = 35641.207
$ws.Cells.Item(100, 11).Value = 35641.207
$ws.Cells.Item(100, 13).Value = -35100.207
$ws.Cells.Item(127, 8).Value = 402039.6
$ws.Cells.Item(127, 9).Value = 402039.6
$ws.Cells.Item(127, 10).Value = 0
$ws.Cells.Item(127, 11).Value = 1206118.8
$ws.Cells.Item(127, 12).Value = 0
$ws.Cells.Item(127, 13).Value = -1201158.8
$ws.Cells.Item(127, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 1842.24
$ws.Cells.Item(132, 9).Value = 1526.7142
$ws.Cells.Item(132, 11).Value = 4580.142599999999
$ws.Cells.Item(132, 13).Value = -2050.142599999999
$ws.Cells.Item(138, 8).Value = 2169.6086
$ws.Cells.Item(138, 9).Value = 1843.3784
$ws.Cells.Item(138, 10).Value = 3510.7778
$ws.Cells.Item(138, 11).Value = 5530.135200000001
$ws.Cells.Item(138, 12).Value = 10532.3334
$ws.Cells.Item(138, 13).Value = -390.1352000000006
$ws.Cells.Item(138, 14).Value = -20812.3334
$ws.Cells.Item(141, 8).Value = 2390.1667
$ws.Cells.Item(141, 9).Value = 2426.5186
$ws.Cells.Item(141, 11).Value = 7279.5558
$ws.Cells.Item(141, 13).Value = -2099.5558

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 5323
$ws.Cells.Item(2, 9).Value = 4021.88
$ws.Cells.Item(2, 10).Value = 16165.667
$ws.Cells.Item(2, 11).Value = 4021.88
$ws.Cells.Item(2, 12).Value = 16165.667
$ws.Cells.Item(2, 13).Value = -3908.88
$ws.Cells.Item(2, 14).Value = -16391.667
$ws.Cells.Item(5, 8).Value = 256
$ws.Cells.Item(5, 9).Value = 232.5
$ws.Cells.Item(5, 10).Value = 350
$ws.Cells.Item(5, 11).Value = 232.5
$ws.Cells.Item(5, 12).Value = 350
$ws.Cells.Item(5, 13).Value = -120.5
$ws.Cells.Item(5, 14).Value = -574
$ws.Cells.Item(32, 8).Value = 2831.3333
$ws.Cells.Item(32, 9).Value = 1555.4615
$ws.Cells.Item(32, 11).Value = 1555.4615
$ws.Cells.Item(32, 13).Value = -1268.4615
$ws.Cells.Item(54, 8).Value = 21332.334
$ws.Cells.Item(54, 10).Value = 21332.334
$ws.Cells.Item(54, 12).Value = 21332.334
$ws.Cells.Item(54, 14).Value = -22870.334
$ws.Cells.Item(58, 8).Value = 36499.5
$ws.Cells.Item(58, 10).Value = 36499.5
$ws.Cells.Item(58, 12).Value = 36499.5
$ws.Cells.Item(58, 14).Value = -37359.5
$ws.Cells.Item(61, 8).Value = 5041.75
$ws.Cells.Item(61, 9).Value = 3207.5454
$ws.Cells.Item(61, 10).Value = 13688.714
$ws.Cells.Item(61, 11).Value = 3207.5454
$ws.Cells.Item(61, 12).Value = 13688.714
$ws.Cells.Item(61, 13).Value = -2995.5454
$ws.Cells.Item(61, 14).Value = -14112.714
$ws.Cells.Item(63, 8).Value = 3833.95
$ws.Cells.Item(63, 9).Value = 3325.0881
$ws.Cells.Item(63, 10).Value = 6717.5
$ws.Cells.Item(63, 11).Value = 3325.0881
$ws.Cells.Item(63, 12).Value = 6717.5
$ws.Cells.Item(63, 13).Value = -2639.0881
$ws.Cells.Item(63, 14).Value = -8089.5
$ws.Cells.Item(66, 8).Value = 3833.95
$ws.Cells.Item(66, 9).Value = 3325.0881
$ws.Cells.Item(66, 10).Value = 6717.5
$ws.Cells.Item(66, 11).Value = 16625.4405
$ws.Cells.Item(66, 12).Value = 33587.5
$ws.Cells.Item(66, 13).Value = -13193.4405
$ws.Cells.Item(66, 14).Value = -40451.5
$ws.Cells.Item(74, 8).Value = 3342.2163
$ws.Cells.Item(74, 9).Value = 2656.1667
$ws.Cells.Item(74, 10).Value = 4608.769
$ws.Cells.Item(74, 11).Value = 2656.1667
$ws.Cells.Item(74, 12).Value = 4608.769
$ws.Cells.Item(74, 13).Value = -1782.1667
$ws.Cells.Item(74, 14).Value = -6356.769
$ws.Cells.Item(77, 8).Value = 3342.2163
$ws.Cells.Item(77, 9).Value = 2656.1667
$ws.Cells.Item(77, 10).Value = 4608.769
$ws.Cells.Item(77, 11).Value = 13280.8335
$ws.Cells.Item(77, 12).Value = 23043.845
$ws.Cells.Item(77, 13).Value = -8912.8335
$ws.Cells.Item(77, 14).Value = -31779.845
$ws.Cells.Item(88, 8).Value = 2595.6667
$ws.Cells.Item(88, 9).Value = 1897.8
$ws.Cells.Item(88, 10).Value = 3468
$ws.Cells.Item(88, 11).Value = 1897.8
$ws.Cells.Item(88, 12).Value = 3468
$ws.Cells.Item(88, 13).Value = -1491.8
$ws.Cells.Item(88, 14).Value = -4280
$ws.Cells.Item(91, 8).Value = 2595.6667
$ws.Cells.Item(91, 9).Value = 1897.8
$ws.Cells.Item(91, 10).Value = 3468
$ws.Cells.Item(91, 11).Value = 1897.8
$ws.Cells.Item(91, 12).Value = 3468
$ws.Cells.Item(91, 13).Value = -493.8
$ws.Cells.Item(91, 14).Value = -6276
$ws.Cells.Item(116, 8).Value = 5323
$ws.Cells.Item(116, 9).Value = 4021.88
$ws.Cells.Item(116, 10).Value = 16165.667
$ws.Cells.Item(116, 11).Value = 4021.88
$ws.Cells.Item(116, 12).Value = 16165.667
$ws.Cells.Item(116, 13).Value = -1727.88
$ws.Cells.Item(116, 14).Value = -20753.667
$ws.Cells.Item(132, 8).Value = 4725.85
$ws.Cells.Item(132, 9).Value = 4725.85
$ws.Cells.Item(132, 11).Value = 14177.55
$ws.Cells.Item(132, 13).Value = -11647.55
$ws.Cells.Item(133, 8).Value = 74999.75
$ws.Cells.Item(133, 10).Value = 74999.75
$ws.Cells.Item(133, 12).Value = 74999.75
$ws.Cells.Item(133, 14).Value = -80059.75
$ws.Cells.Item(136, 8).Value = 5041.75
$ws.Cells.Item(136, 9).Value = 3207.5454
$ws.Cells.Item(136, 10).Value = 13688.714
$ws.Cells.Item(136, 11).Value = 9622.6362
$ws.Cells.Item(136, 12).Value = 41066.142
$ws.Cells.Item(136, 13).Value = -7072.636200000001
$ws.Cells.Item(136, 14).Value = -46166.142

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 5323
$ws.Cells.Item(3, 9).Value = 4021.88
$ws.Cells.Item(3, 10).Value = 16165.667
$ws.Cells.Item(3, 11).Value = 4021.88
$ws.Cells.Item(3, 12).Value = 16165.667
$ws.Cells.Item(3, 13).Value = -3907.88
$ws.Cells.Item(3, 14).Value = -16393.667
$ws.Cells.Item(4, 8).Value = 256
$ws.Cells.Item(4, 9).Value = 232.5
$ws.Cells.Item(4, 10).Value = 350
$ws.Cells.Item(4, 11).Value = 232.5
$ws.Cells.Item(4, 12).Value = 350
$ws.Cells.Item(4, 13).Value = -117.5
$ws.Cells.Item(4, 14).Value = -580
$ws.Cells.Item(54, 8).Value = 2349.75
$ws.Cells.Item(54, 10).Value = 0
$ws.Cells.Item(54, 12).Value = 0
$ws.Cells.Item(54, 14).ClearContents()
$ws.Cells.Item(82, 8).Value = 28321.646
$ws.Cells.Item(82, 9).Value = 15237.333
$ws.Cells.Item(82, 10).Value = 31125.428
$ws.Cells.Item(82, 11).Value = 15237.333
$ws.Cells.Item(82, 12).Value = 31125.428
$ws.Cells.Item(82, 13).Value = -14854.333
$ws.Cells.Item(82, 14).Value = -31891.428
$ws.Cells.Item(85, 8).Value = 28321.646
$ws.Cells.Item(85, 9).Value = 15237.333
$ws.Cells.Item(85, 10).Value = 31125.428
$ws.Cells.Item(85, 11).Value = 15237.333
$ws.Cells.Item(85, 12).Value = 31125.428
$ws.Cells.Item(85, 13).Value = -13911.333
$ws.Cells.Item(85, 14).Value = -33777.428
$ws.Cells.Item(94, 8).Value = 652.35596
$ws.Cells.Item(94, 9).Value = 661.2157
$ws.Cells.Item(94, 10).Value = 595.875
$ws.Cells.Item(94, 11).Value = 661.2157
$ws.Cells.Item(94, 12).Value = 595.875
$ws.Cells.Item(94, 13).Value = -210.2157
$ws.Cells.Item(94, 14).Value = -1497.875
$ws.Cells.Item(99, 8).Value = 4424.0713
$ws.Cells.Item(99, 9).Value = 3504.6667
$ws.Cells.Item(99, 11).Value = 3504.6667
$ws.Cells.Item(99, 13).Value = -2006.6667
$ws.Cells.Item(105, 8).Value = 2816.291
$ws.Cells.Item(105, 9).Value = 3330.2173
$ws.Cells.Item(105, 10).Value = 2446.9062
$ws.Cells.Item(105, 11).Value = 3330.2173
$ws.Cells.Item(105, 12).Value = 2446.9062
$ws.Cells.Item(105, 13).Value = -1583.2173
$ws.Cells.Item(105, 14).Value = -5940.906199999999
$ws.Cells.Item(107, 8).Value = 2246.3333
$ws.Cells.Item(107, 9).Value = 1662.3334
$ws.Cells.Item(107, 11).Value = 1662.3334
$ws.Cells.Item(107, 13).Value = 257.6666
$ws.Cells.Item(112, 8).Value = 0
$ws.Cells.Item(112, 10).Value = 0
$ws.Cells.Item(112, 12).Value = 0
$ws.Cells.Item(112, 14).ClearContents()
$ws.Cells.Item(134, 8).Value = 11062.275
$ws.Cells.Item(134, 9).Value = 5834.3438
$ws.Cells.Item(134, 10).Value = 31974
$ws.Cells.Item(134, 11).Value = 17503.0314
$ws.Cells.Item(134, 12).Value = 95922
$ws.Cells.Item(134, 13).Value = -14968.0314
$ws.Cells.Item(134, 14).Value = -100992

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(5, 8).Value = 405.6
$ws.Cells.Item(5, 9).Value = 233.33333
$ws.Cells.Item(5, 10).Value = 479.42856
$ws.Cells.Item(5, 11).Value = 233.33333
$ws.Cells.Item(5, 12).Value = 479.42856
$ws.Cells.Item(5, 13).Value = -121.33333
$ws.Cells.Item(5, 14).Value = -703.4285600000001
$ws.Cells.Item(7, 8).Value = 194
$ws.Cells.Item(7, 9).Value = 196.5
$ws.Cells.Item(7, 10).Value = 189
$ws.Cells.Item(7, 11).Value = 196.5
$ws.Cells.Item(7, 12).Value = 189
$ws.Cells.Item(7, 13).Value = -83.5
$ws.Cells.Item(7, 14).Value = -415
$ws.Cells.Item(10, 8).Value = 1993.1538
$ws.Cells.Item(10, 9).Value = 2600.75
$ws.Cells.Item(10, 10).Value = 1723.1111
$ws.Cells.Item(10, 11).Value = 2600.75
$ws.Cells.Item(10, 12).Value = 1723.1111
$ws.Cells.Item(10, 13).Value = -2461.75
$ws.Cells.Item(10, 14).Value = -2001.1111
$ws.Cells.Item(16, 8).Value = 992.125
$ws.Cells.Item(16, 9).Value = 696.4
$ws.Cells.Item(16, 10).Value = 1485
$ws.Cells.Item(16, 11).Value = 696.4
$ws.Cells.Item(16, 12).Value = 1485
$ws.Cells.Item(16, 13).Value = -409.4
$ws.Cells.Item(16, 14).Value = -2059
$ws.Cells.Item(22, 8).Value = 186.40909
$ws.Cells.Item(22, 9).Value = 190.76471
$ws.Cells.Item(22, 11).Value = 190.76471
$ws.Cells.Item(22, 13).Value = 159.23529
$ws.Cells.Item(53, 8).Value = 89999
$ws.Cells.Item(53, 10).Value = 89999
$ws.Cells.Item(53, 12).Value = 89999
$ws.Cells.Item(53, 14).Value = -91213
$ws.Cells.Item(58, 8).Value = 6172.8887
$ws.Cells.Item(58, 9).Value = 4213.5
$ws.Cells.Item(58, 10).Value = 21848
$ws.Cells.Item(58, 11).Value = 4213.5
$ws.Cells.Item(58, 12).Value = 21848
$ws.Cells.Item(58, 13).Value = -4010.5
$ws.Cells.Item(58, 14).Value = -22254
$ws.Cells.Item(62, 8).Value = 58075.74
$ws.Cells.Item(62, 9).Value = 94102.73
$ws.Cells.Item(62, 10).Value = 8538.625
$ws.Cells.Item(62, 11).Value = 94102.73
$ws.Cells.Item(62, 12).Value = 8538.625
$ws.Cells.Item(62, 13).Value = -93478.73
$ws.Cells.Item(62, 14).Value = -9786.625
$ws.Cells.Item(65, 8).Value = 58075.74
$ws.Cells.Item(65, 9).Value = 94102.73
$ws.Cells.Item(65, 10).Value = 8538.625
$ws.Cells.Item(65, 11).Value = 470513.65
$ws.Cells.Item(65, 12).Value = 42693.125
$ws.Cells.Item(65, 13).Value = -467393.65
$ws.Cells.Item(65, 14).Value = -48933.125
$ws.Cells.Item(86, 9).Value = 4000
$ws.Cells.Item(86, 11).Value = 4000
$ws.Cells.Item(86, 13).Value = -2877
$ws.Cells.Item(89, 9).Value = 4000
$ws.Cells.Item(89, 11).Value = 20000
$ws.Cells.Item(89, 13).Value = -14384
$ws.Cells.Item(94, 8).Value = 1080.9565
$ws.Cells.Item(94, 9).Value = 1421.5
$ws.Cells.Item(94, 10).Value = 1009.2632
$ws.Cells.Item(94, 11).Value = 1421.5
$ws.Cells.Item(94, 12).Value = 1009.2632
$ws.Cells.Item(94, 13).Value = -970.5
$ws.Cells.Item(94, 14).Value = -1911.2632
$ws.Cells.Item(105, 8).Value = 1987.4667
$ws.Cells.Item(105, 9).Value = 2160.2
$ws.Cells.Item(105, 11).Value = 2160.2
$ws.Cells.Item(105, 13).Value = -413.1999999999998
$ws.Cells.Item(107, 8).Value = 2102.6
$ws.Cells.Item(107, 9).Value = 0
$ws.Cells.Item(107, 10).Value = 2102.6
$ws.Cells.Item(107, 11).Value = 0
$ws.Cells.Item(107, 12).Value = 2102.6
$ws.Cells.Item(107, 13).ClearContents()
$ws.Cells.Item(107, 14).Value = -5942.6
$ws.Cells.Item(113, 8).Value = 992.125
$ws.Cells.Item(113, 9).Value = 696.4
$ws.Cells.Item(113, 10).Value = 1485
$ws.Cells.Item(113, 11).Value = 696.4
$ws.Cells.Item(113, 12).Value = 1485
$ws.Cells.Item(113, 13).Value = 1473.6
$ws.Cells.Item(113, 14).Value = -5825
$ws.Cells.Item(133, 8).Value = 44159
$ws.Cells.Item(133, 10).Value = 44159
$ws.Cells.Item(133, 12).Value = 44159
$ws.Cells.Item(133, 14).Value = -49219
$ws.Cells.Item(134, 8).Value = 6399.72
$ws.Cells.Item(134, 9).Value = 4860.1177
$ws.Cells.Item(134, 10).Value = 9671.375
$ws.Cells.Item(134, 11).Value = 14580.3531
$ws.Cells.Item(134, 12).Value = 29014.125
$ws.Cells.Item(134, 13).Value = -12045.3531
$ws.Cells.Item(134, 14).Value = -34084.125
$ws.Cells.Item(136, 8).Value = 6172.8887
$ws.Cells.Item(136, 9).Value = 4213.5
$ws.Cells.Item(136, 10).Value = 21848
$ws.Cells.Item(136, 11).Value = 12640.5
$ws.Cells.Item(136, 12).Value = 65544
$ws.Cells.Item(136, 13).Value = -10090.5
$ws.Cells.Item(136, 14).Value = -70644

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 181.62962
$ws.Cells.Item(2, 9).Value = 24.333334
$ws.Cells.Item(2, 10).Value = 732.1667
$ws.Cells.Item(2, 11).Value = 146.000004
$ws.Cells.Item(2, 12).Value = 4393.0002
$ws.Cells.Item(2, 13).Value = -33.00000399999999
$ws.Cells.Item(2, 14).Value = -4619.0002
$ws.Cells.Item(7, 8).Value = 595
$ws.Cells.Item(7, 10).Value = 595
$ws.Cells.Item(7, 12).Value = 1785
$ws.Cells.Item(7, 14).Value = -2009
$ws.Cells.Item(34, 8).Value = 2146.4285
$ws.Cells.Item(34, 10).Value = 5971.143
$ws.Cells.Item(34, 12).Value = 17913.429
$ws.Cells.Item(34, 14).Value = -18081.429
$ws.Cells.Item(39, 8).Value = 7800
$ws.Cells.Item(39, 10).Value = 8990
$ws.Cells.Item(39, 12).Value = 26970
$ws.Cells.Item(39, 14).Value = -27558
$ws.Cells.Item(107, 8).Value = 676.9524
$ws.Cells.Item(107, 10).Value = 817.6667
$ws.Cells.Item(107, 12).Value = 2453.0001
$ws.Cells.Item(107, 14).Value = -6293.0001
$ws.Cells.Item(121, 8).Value = 1887.9231
$ws.Cells.Item(121, 9).Value = 1573.2
$ws.Cells.Item(121, 10).Value = 2937
$ws.Cells.Item(121, 11).Value = 4719.6
$ws.Cells.Item(121, 12).Value = 8811
$ws.Cells.Item(121, 13).Value = -3409.6
$ws.Cells.Item(121, 14).Value = -11431
$ws.Cells.Item(122, 8).Value = 10001396
$ws.Cells.Item(122, 10).Value = 14286855
$ws.Cells.Item(122, 12).Value = 128581695
$ws.Cells.Item(122, 14).Value = -128586595
$ws.Cells.Item(137, 8).Value = 5355.5
$ws.Cells.Item(137, 9).Value = 4371.8
$ws.Cells.Item(137, 11).Value = 13115.4
$ws.Cells.Item(137, 13).Value = -8015.400000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(69, 8).Value = 60000
$ws.Cells.Item(69, 10).Value = 60000
$ws.Cells.Item(69, 12).Value = 60000
$ws.Cells.Item(69, 14).Value = -61498
$ws.Cells.Item(72, 8).Value = 60000
$ws.Cells.Item(72, 10).Value = 60000
$ws.Cells.Item(72, 12).Value = 180000
$ws.Cells.Item(72, 14).Value = -187488
$ws.Cells.Item(80, 8).Value = 8765.526
$ws.Cells.Item(80, 9).Value = 7131.7856
$ws.Cells.Item(80, 10).Value = 13340
$ws.Cells.Item(80, 11).Value = 7131.7856
$ws.Cells.Item(80, 12).Value = 13340
$ws.Cells.Item(80, 13).Value = -6133.7856
$ws.Cells.Item(80, 14).Value = -15336
$ws.Cells.Item(83, 8).Value = 8765.526
$ws.Cells.Item(83, 9).Value = 7131.7856
$ws.Cells.Item(83, 10).Value = 13340
$ws.Cells.Item(83, 11).Value = 35658.928
$ws.Cells.Item(83, 12).Value = 66700
$ws.Cells.Item(83, 13).Value = -30666.928
$ws.Cells.Item(83, 14).Value = -76684
$ws.Cells.Item(97, 8).Value = 1117.931
$ws.Cells.Item(97, 9).Value = 943.4762
$ws.Cells.Item(97, 10).Value = 1575.875
$ws.Cells.Item(97, 11).Value = 943.4762
$ws.Cells.Item(97, 12).Value = 1575.875
$ws.Cells.Item(97, 13).Value = -447.4761999999999
$ws.Cells.Item(97, 14).Value = -2567.875
$ws.Cells.Item(122, 8).Value = 2404.3333
$ws.Cells.Item(122, 9).Value = 1435
$ws.Cells.Item(122, 11).Value = 4305
$ws.Cells.Item(122, 13).Value = -1855
$ws.Cells.Item(132, 8).Value = 10709.962
$ws.Cells.Item(132, 9).Value = 11863.521
$ws.Cells.Item(132, 10).Value = 1866
$ws.Cells.Item(132, 11).Value = 35590.563
$ws.Cells.Item(132, 12).Value = 5598
$ws.Cells.Item(132, 13).Value = -33060.563
$ws.Cells.Item(132, 14).Value = -10658

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 11029.235
$ws.Cells.Item(7, 9).Value = 11924.333
$ws.Cells.Item(7, 10).Value = 8881
$ws.Cells.Item(7, 11).Value = 11924.333
$ws.Cells.Item(7, 12).Value = 8881
$ws.Cells.Item(7, 13).Value = -11812.333
$ws.Cells.Item(7, 14).Value = -9105
$ws.Cells.Item(16, 8).Value = 6821
$ws.Cells.Item(16, 9).Value = 709.7857
$ws.Cells.Item(16, 10).Value = 49599.5
$ws.Cells.Item(16, 11).Value = 709.7857
$ws.Cells.Item(16, 12).Value = 49599.5
$ws.Cells.Item(16, 13).Value = -539.7857
$ws.Cells.Item(16, 14).Value = -49939.5
$ws.Cells.Item(22, 8).Value = 16483.5
$ws.Cells.Item(22, 9).Value = 90000
$ws.Cells.Item(22, 10).Value = 1780.2
$ws.Cells.Item(22, 11).Value = 90000
$ws.Cells.Item(22, 12).Value = 1780.2
$ws.Cells.Item(22, 13).Value = -89705
$ws.Cells.Item(22, 14).Value = -2370.2
$ws.Cells.Item(27, 8).Value = 16483.5
$ws.Cells.Item(27, 9).Value = 90000
$ws.Cells.Item(27, 10).Value = 1780.2
$ws.Cells.Item(27, 11).Value = 90000
$ws.Cells.Item(27, 12).Value = 1780.2
$ws.Cells.Item(27, 13).Value = -89893
$ws.Cells.Item(27, 14).Value = -1994.2
$ws.Cells.Item(46, 8).Value = 1359.7812
$ws.Cells.Item(46, 9).Value = 977.7143
$ws.Cells.Item(46, 10).Value = 1466.76
$ws.Cells.Item(46, 11).Value = 977.7143
$ws.Cells.Item(46, 12).Value = 1466.76
$ws.Cells.Item(46, 13).Value = -789.7143
$ws.Cells.Item(46, 14).Value = -1842.76
$ws.Cells.Item(55, 8).Value = 150.05263
$ws.Cells.Item(55, 9).Value = 99.27273
$ws.Cells.Item(55, 10).Value = 219.875
$ws.Cells.Item(55, 11).Value = 99.27273
$ws.Cells.Item(55, 12).Value = 219.875
$ws.Cells.Item(55, 13).Value = 73.72727
$ws.Cells.Item(55, 14).Value = -565.875
$ws.Cells.Item(68, 8).Value = 2580.3818
$ws.Cells.Item(68, 9).Value = 2356.9783
$ws.Cells.Item(68, 10).Value = 3722.2222
$ws.Cells.Item(68, 11).Value = 2356.9783
$ws.Cells.Item(68, 12).Value = 3722.2222
$ws.Cells.Item(68, 13).Value = -1607.9783
$ws.Cells.Item(68, 14).Value = -5220.2222
$ws.Cells.Item(71, 8).Value = 2580.3818
$ws.Cells.Item(71, 9).Value = 2356.9783
$ws.Cells.Item(71, 10).Value = 3722.2222
$ws.Cells.Item(71, 11).Value = 11784.8915
$ws.Cells.Item(71, 12).Value = 18611.111
$ws.Cells.Item(71, 13).Value = -8040.891500000001
$ws.Cells.Item(71, 14).Value = -26099.111
$ws.Cells.Item(82, 8).Value = 2034.619
$ws.Cells.Item(82, 9).Value = 1569.6666
$ws.Cells.Item(82, 10).Value = 2654.5557
$ws.Cells.Item(82, 11).Value = 1569.6666
$ws.Cells.Item(82, 12).Value = 2654.5557
$ws.Cells.Item(82, 13).Value = -1208.6666
$ws.Cells.Item(82, 14).Value = -3376.5557
$ws.Cells.Item(85, 8).Value = 2034.619
$ws.Cells.Item(85, 9).Value = 1569.6666
$ws.Cells.Item(85, 10).Value = 2654.5557
$ws.Cells.Item(85, 11).Value = 1569.6666
$ws.Cells.Item(85, 12).Value = 2654.5557
$ws.Cells.Item(85, 13).Value = -321.6666
$ws.Cells.Item(85, 14).Value = -5150.5557
$ws.Cells.Item(93, 8).Value = 3641
$ws.Cells.Item(93, 9).Value = 3711.875
$ws.Cells.Item(93, 10).Value = 3479
$ws.Cells.Item(93, 11).Value = 3711.875
$ws.Cells.Item(93, 12).Value = 3479
$ws.Cells.Item(93, 13).Value = -2463.875
$ws.Cells.Item(93, 14).Value = -5975
$ws.Cells.Item(100, 8).Value = 8621.546
$ws.Cells.Item(100, 9).Value = 6927.4443
$ws.Cells.Item(100, 11).Value = 6927.4443
$ws.Cells.Item(100, 13).Value = -6386.4443
$ws.Cells.Item(122, 8).Value = 1447.5
$ws.Cells.Item(122, 9).Value = 1263.3334
$ws.Cells.Item(122, 11).Value = 3790.0002
$ws.Cells.Item(122, 13).Value = -1340.0002
$ws.Cells.Item(126, 8).Value = 11029.235
$ws.Cells.Item(126, 9).Value = 11924.333
$ws.Cells.Item(126, 10).Value = 8881
$ws.Cells.Item(126, 11).Value = 35772.999
$ws.Cells.Item(126, 12).Value = 26643
$ws.Cells.Item(126, 13).Value = -33302.999
$ws.Cells.Item(126, 14).Value = -31583
$ws.Cells.Item(132, 8).Value = 4184.65
$ws.Cells.Item(132, 9).Value = 3890.8333
$ws.Cells.Item(132, 11).Value = 11672.4999
$ws.Cells.Item(132, 13).Value = -9142.499899999999
$ws.Cells.Item(136, 8).Value = 2726.3542
$ws.Cells.Item(136, 9).Value = 2407.439
$ws.Cells.Item(136, 11).Value = 7222.316999999999
$ws.Cells.Item(136, 13).Value = -4672.316999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(74, 8).Value = 26264.572
$ws.Cells.Item(74, 10).Value = 27325.334
$ws.Cells.Item(74, 12).Value = 27325.334
$ws.Cells.Item(74, 14).Value = -29197.334
$ws.Cells.Item(75, 8).Value = 38789.5
$ws.Cells.Item(75, 10).Value = 26629.5
$ws.Cells.Item(75, 12).Value = 26629.5
$ws.Cells.Item(75, 14).Value = -28501.5
$ws.Cells.Item(77, 8).Value = 26264.572
$ws.Cells.Item(77, 10).Value = 27325.334
$ws.Cells.Item(77, 12).Value = 81976.002
$ws.Cells.Item(77, 14).Value = -91336.002
$ws.Cells.Item(78, 8).Value = 38789.5
$ws.Cells.Item(78, 10).Value = 26629.5
$ws.Cells.Item(78, 12).Value = 79888.5
$ws.Cells.Item(78, 14).Value = -89248.5
$ws.Cells.Item(81, 8).Value = 2497.6428
$ws.Cells.Item(81, 9).Value = 2251.889
$ws.Cells.Item(81, 10).Value = 2940
$ws.Cells.Item(81, 11).Value = 4503.778
$ws.Cells.Item(81, 12).Value = 5880
$ws.Cells.Item(81, 13).Value = -3442.778
$ws.Cells.Item(81, 14).Value = -8002
$ws.Cells.Item(84, 8).Value = 2497.6428
$ws.Cells.Item(84, 9).Value = 2251.889
$ws.Cells.Item(84, 10).Value = 2940
$ws.Cells.Item(84, 11).Value = 22518.89
$ws.Cells.Item(84, 12).Value = 29400
$ws.Cells.Item(84, 13).Value = -17214.89
$ws.Cells.Item(84, 14).Value = -40008
$ws.Cells.Item(126, 8).Value = 6974.4546
$ws.Cells.Item(126, 9).Value = 6974.4546
$ws.Cells.Item(126, 10).Value = 0
$ws.Cells.Item(126, 11).Value = 20923.3638
$ws.Cells.Item(126, 12).Value = 0
$ws.Cells.Item(126, 13).Value = -18453.3638
$ws.Cells.Item(126, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 15125.055
$ws.Cells.Item(132, 9).Value = 7419.346
$ws.Cells.Item(132, 11).Value = 22258.038
$ws.Cells.Item(132, 13).Value = -19728.038
$ws.Cells.Item(136, 8).Value = 869.7632
$ws.Cells.Item(136, 9).Value = 882.9167
$ws.Cells.Item(136, 11).Value = 2648.7501
$ws.Cells.Item(136, 13).Value = -98.7501000000002
$ws.Cells.Item(137, 8).Value = 96046.664
$ws.Cells.Item(137, 9).Value = 0
$ws.Cells.Item(137, 10).Value = 96046.664
$ws.Cells.Item(137, 11).Value = 0
$ws.Cells.Item(137, 12).Value = 96046.664
$ws.Cells.Item(137, 13).ClearContents()
$ws.Cells.Item(137, 14).Value = -106246.664
